$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update values that changed on existing rows ---
# Row 2: quantidade_atipica / desvio_padrao recalculated
$ws.Range("G2").Value = -40
$ws.Range("I2").Value = 0.14

# Row 3: quantidade_atipica recalculated
$ws.Range("G3").Value = -112

# --- Append new row 5 (new atypical sale record for 2025-06-12) ---
# A5 and D5 must stay text (date-like / numeric-looking codes), so force
# text format before assigning so Excel doesn't auto-convert them.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"

$ws.Range("A5").Value = "2025-06-12"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "BEMOL S/A"
$ws.Range("D5").Value = "358537"
$ws.Range("E5").Value = 13588
$ws.Range("F5").Value = "CANETA STYLLUS ACTIVA AGOLD"
$ws.Range("G5").Value = -4
$ws.Range("H5").Value = 1.05
$ws.Range("I5").Value = 0.21
